$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 443 (pushes the former rows 443:465 down to 444:466,
# and extends the sheet dimension from A1:R465 to A1:R466).
$ws.Rows.Item(443).Insert()

# Populate the new row 443 with the new data record. Columns A, B, C, E, F, G, N, Q, R
# carry the same constant values as the surrounding "Vega Modelo de Temuco / Zapallo" rows.
$ws.Range("A443").Value = 10
$ws.Range("B443").Value = "Vega Modelo de Temuco"
$ws.Range("C443").Value = "La Araucanía"
$ws.Range("D443").Value = 44585
$ws.Range("E443").Value = 9
$ws.Range("F443").Value = 100112045
$ws.Range("G443").Value = "Zapallo"
$ws.Range("H443").Value = "Paine"
$ws.Range("I443").Value = "1a nueva(o)"
$ws.Range("J443").Value = 650
$ws.Range("K443").Value = 300
$ws.Range("L443").Value = 300
$ws.Range("M443").Value = 300
$ws.Range("N443").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O443").Value = "Región de O'Higgins"
$ws.Range("P443").Value = 300
$ws.Range("Q443").Value = 1
$ws.Range("R443").Value = "Hortaliza"
